$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.634144425392151
$ws.Range("B1").Value = 4.041515827178955
$ws.Range("C1").Value = 3.865645885467529
$ws.Range("D1").Value = 1.64033579826355
$ws.Range("E1").Value = 1.062463521957397
